# Slide 5 ("System Architecture") - "Content Placeholder 2" bullet list.
# The last two bullets:
#   "Database: MongoDB for storing users, expenses, and income"
#   "Simple architecture diagram (placeholder)"
# are replaced by a single bullet:
#   "Database: MongoDB for storing users, expenses, and income"
# (i.e. the placeholder bullet is dropped), and that bullet's text ends
# up split across two runs ("...expenses, " / "and income") to match the
# authored edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame.TextRange

$tr.Text = "Frontend: Handles UI, user interactions, and API requests`r" + `
           "Backend: REST API, authentication, business logic`r" + `
           "Database: MongoDB for storing users, expenses, and income"

# Split the third paragraph's single run into two runs, breaking right
# before "and income", matching the authored diff.
$para3 = $tr.Paragraphs(3, 1)
$splitAt = $para3.Text.IndexOf("and income") + 1
$firstRun = $para3.Characters(1, $splitAt - 1)
$firstRun.Text = "Database: MongoDB for storing users, expenses, "
